$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1977715877437326
$ws.Range("C2").Value = 0.5459610027855153
$ws.Range("J2").Value = 0.005571030640668524
$ws.Range("P2").Value = 0.1364902506963788
$ws.Range("S2").Value = 0.1142061281337047
$ws.Range("C3").Value = 0.01015228426395939
$ws.Range("J3").Value = 0.03045685279187817
$ws.Range("P3").Value = 0.7411167512690355
$ws.Range("S3").Value = 0.2182741116751269
$ws.Range("J4").Value = 0.07272727272727272
$ws.Range("P4").Value = 0.6181818181818182
$ws.Range("S4").Value = 0.3090909090909091
$ws.Range("B6").Value = 0.06751054852320675
$ws.Range("D6").Value = 0.02953586497890295
$ws.Range("F6").Value = 0.05063291139240506
$ws.Range("J6").Value = 0.2911392405063291
$ws.Range("O6").Value = 0.02531645569620253
$ws.Range("Q6").Value = 0.1856540084388186
$ws.Range("R6").Value = 0.0759493670886076
$ws.Range("S6").Value = 0.2742616033755274
$ws.Range("B7").Value = 0.1142857142857143
$ws.Range("D7").Value = 0.02857142857142857
$ws.Range("F7").Value = 0.04285714285714286
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("O7").Value = 0.04285714285714286
$ws.Range("Q7").Value = 0.1142857142857143
$ws.Range("R7").Value = 0.05714285714285714
$ws.Range("S7").Value = 0.4333333333333333
$ws.Range("B8").Value = 0.09664694280078895
$ws.Range("D8").Value = 0.01183431952662722
$ws.Range("F8").Value = 0.04142011834319527
$ws.Range("J8").Value = 0.1479289940828402
$ws.Range("O8").Value = 0.02564102564102564
$ws.Range("Q8").Value = 0.1637080867850099
$ws.Range("R8").Value = 0.08678500986193294
$ws.Range("S8").Value = 0.4260355029585799
$ws.Range("B9").Value = 0.0847457627118644
$ws.Range("D9").Value = 0.01694915254237288
$ws.Range("F9").Value = 0.0635593220338983
$ws.Range("J9").Value = 0.1694915254237288
$ws.Range("O9").Value = 0.0211864406779661
$ws.Range("Q9").Value = 0.1694915254237288
$ws.Range("R9").Value = 0.08050847457627118
$ws.Range("S9").Value = 0.3940677966101695
$ws.Range("B10").Value = 0.1203319502074689
$ws.Range("D10").Value = 0.02282157676348548
$ws.Range("E10").Value = 0.001383125864453665
$ws.Range("F10").Value = 0.0656984785615491
$ws.Range("J10").Value = 0.1334716459197787
$ws.Range("O10").Value = 0.01867219917012448
$ws.Range("Q10").Value = 0.2019363762102351
$ws.Range("R10").Value = 0.07330567081604426
$ws.Range("S10").Value = 0.3623789764868603
$ws.Range("G11").Value = 0.1529411764705882
$ws.Range("J11").Value = 0.06764705882352941
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.5647058823529412
$ws.Range("S11").Value = 0.01470588235294118
$ws.Range("G12").Value = 0.6811594202898551
$ws.Range("J12").Value = 0.2318840579710145
$ws.Range("K12").Value = 0.004830917874396135
$ws.Range("L12").Value = 0.02415458937198068
$ws.Range("S12").Value = 0.05797101449275362
$ws.Range("G13").Value = 0.65
$ws.Range("J13").Value = 0.35
$ws.Range("F15").Value = 0.02788844621513944
$ws.Range("H15").Value = 0.1434262948207171
$ws.Range("I15").Value = 0.06772908366533864
$ws.Range("J15").Value = 0.3266932270916335
$ws.Range("K15").Value = 0.05577689243027888
$ws.Range("M15").Value = 0.02390438247011952
$ws.Range("O15").Value = 0.06374501992031872
$ws.Range("S15").Value = 0.2908366533864542
$ws.Range("F16").Value = 0.01785714285714286
$ws.Range("H16").Value = 0.1428571428571428
$ws.Range("I16").Value = 0.1160714285714286
$ws.Range("J16").Value = 0.4241071428571428
$ws.Range("K16").Value = 0.08928571428571429
$ws.Range("M16").Value = 0.004464285714285714
$ws.Range("O16").Value = 0.05357142857142857
$ws.Range("S16").Value = 0.1517857142857143
$ws.Range("F17").Value = 0.02928870292887029
$ws.Range("H17").Value = 0.1715481171548117
$ws.Range("I17").Value = 0.1192468619246862
$ws.Range("J17").Value = 0.3849372384937239
$ws.Range("K17").Value = 0.1066945606694561
$ws.Range("M17").Value = 0.008368200836820083
$ws.Range("O17").Value = 0.07112970711297072
$ws.Range("S17").Value = 0.1087866108786611
$ws.Range("F18").Value = 0.01522842639593909
$ws.Range("H18").Value = 0.1573604060913706
$ws.Range("I18").Value = 0.08121827411167512
$ws.Range("J18").Value = 0.4162436548223351
$ws.Range("K18").Value = 0.116751269035533
$ws.Range("M18").Value = 0.01015228426395939
$ws.Range("O18").Value = 0.07106598984771574
$ws.Range("S18").Value = 0.1319796954314721
$ws.Range("F19").Value = 0.02083333333333333
$ws.Range("H19").Value = 0.2305555555555556
$ws.Range("I19").Value = 0.08333333333333333
$ws.Range("J19").Value = 0.3569444444444445
$ws.Range("K19").Value = 0.1097222222222222
$ws.Range("M19").Value = 0.01875
$ws.Range("N19").Value = 0.001388888888888889
$ws.Range("O19").Value = 0.06041666666666667
$ws.Range("S19").Value = 0.1180555555555556
